$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (style index 0 / default, no custom formatting) used to
# strip the "quote prefix" formatting flag that Excel applies when a
# numeric-looking string is typed into a cell, so cells keep matching the
# plain default style of the original workbook.
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "'28.105.78"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "'1.799.97"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'311.18"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.5108"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  -2.60%  "
$ws.Range("D8").Value = "'0.3890"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("D9").Value = "'0.07733"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "'40.92"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("D12").Value = "'6.313"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "'20.24"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "'1.795.71"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").Value = "'7.260"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "'92.04"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").Value = "'0.06577"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "'17.20"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "'5.961"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'28.129.35"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "'11.05"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("D25").Value = "'2.254"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").Value = "'160.58"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").Value = "'2.417"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +2.08%  "
$ws.Range("D28").Value = "'2.005.33"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").Value = "'20.23"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").Value = "'127.13"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +2.79%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").Value = "'1.046"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  -1.40%  "
$ws.Range("D33").Value = "'3.645"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").Value = "'0.07025"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("D36").Value = "'9.031"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +4.05%  "
$ws.Range("D37").Value = "'0.02340"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").Value = "'0.2154"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.008"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'11.46"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  -5.63%  "
$ws.Range("D41").Value = "'0.6109"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("D42").Value = "'1.002"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").Value = "'1.151"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").Value = "'13.07"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").Value = "'1.297"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  -6.23%  "
$ws.Range("D46").Value = "'0.5900"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  -2.15%  "
$ws.Range("D47").Value = "'3.719"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").Value = "'125.03"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").Value = "'1.198"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "'1.899"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("D51").Value = "'0.06734"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -1.31%  "
